# Update Wild Card round simulated stats for Rams 2021 Target Depth Data
$wb = $excel.ActiveWorkbook

# OFF sheet, row 2 ("H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 240
$wsOff.Range("C2").Value = 177
$wsOff.Range("D2").Value = 71
$wsOff.Range("E2").Value = 37
$wsOff.Range("F2").Value = 2
$wsOff.Range("G2").Value = 5

# DEF sheet, row 2 ("H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 243
$wsDef.Range("C2").Value = 165
$wsDef.Range("D2").Value = 64
$wsDef.Range("E2").Value = 24
$wsDef.Range("F2").Value = 7
$wsDef.Range("G2").Value = 9
